# Update the "Förändrad" (changed) date column C for every data row, and
# append the Beteckning (column A) value as a second argument to every
# HYPERLINK() formula found in columns S-Y, matching the new link-text
# behaviour introduced upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

$hyperlinkCols = @("S", "T", "U", "V", "W", "X", "Y")

for ($r = 2; $r -le $lastRow; $r++) {
    # Column C: new "changed" timestamp (serial date 45184 -> 45186).
    $ws.Cells.Item($r, 3).Value = 45186

    $beteckning = $ws.Cells.Item($r, 1).Value2
    if ($beteckning -eq $null -or $beteckning -eq "") {
        continue
    }

    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Range("$col$r")
        $f = $cell.Formula
        if ($f -eq $null -or $f -eq "") {
            continue
        }
        if ($f.StartsWith("=HYPERLINK(") -and -not $f.Contains(",")) {
            $newFormula = $f.Substring(0, $f.Length - 1) + ', "' + $beteckning + '")'
            $cell.Formula = $newFormula
        }
    }
}
